$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H62").Value = 34362.906
$ws.Range("I62").Value = 52088.65
$ws.Range("K62").Value = 52088.65
$ws.Range("M62").Value = -51464.65
$ws.Range("H65").Value = 34362.906
$ws.Range("I65").Value = 52088.65
$ws.Range("K65").Value = 260443.25
$ws.Range("M65").Value = -257323.25
$ws.Range("H80").Value = 299629.97
$ws.Range("I80").Value = 559.06665
$ws.Range("J80").Value = 707453.9399999999
$ws.Range("K80").Value = 1677.19995
$ws.Range("L80").Value = 2122361.82
$ws.Range("M80").Value = -679.1999499999999
$ws.Range("N80").Value = -2124357.82
$ws.Range("H83").Value = 299629.97
$ws.Range("I83").Value = 559.06665
$ws.Range("J83").Value = 707453.9399999999
$ws.Range("K83").Value = 5031.59985
$ws.Range("L83").Value = 6367085.459999999
$ws.Range("M83").Value = -39.59984999999961
$ws.Range("N83").Value = -6377069.459999999
$ws.Range("H86").Value = 50003616
$ws.Range("I86").Value = 3272.8572
$ws.Range("J86").Value = 76926880
$ws.Range("K86").Value = 3272.8572
$ws.Range("L86").Value = 76926880
$ws.Range("M86").Value = -2149.8572
$ws.Range("N86").Value = -76929126
$ws.Range("H88").Value = 1589203.8
$ws.Range("I88").Value = 3002.8
$ws.Range("J88").Value = 3175404.8
$ws.Range("K88").Value = 3002.8
$ws.Range("L88").Value = 3175404.8
$ws.Range("M88").Value = -2596.8
$ws.Range("N88").Value = -3176216.8
$ws.Range("H89").Value = 50003616
$ws.Range("I89").Value = 3272.8572
$ws.Range("J89").Value = 76926880
$ws.Range("K89").Value = 16364.286
$ws.Range("L89").Value = 384634400
$ws.Range("M89").Value = -10748.286
$ws.Range("N89").Value = -384645632
$ws.Range("H91").Value = 1589203.8
$ws.Range("I91").Value = 3002.8
$ws.Range("J91").Value = 3175404.8
$ws.Range("K91").Value = 3002.8
$ws.Range("L91").Value = 3175404.8
$ws.Range("M91").Value = -1598.8
$ws.Range("N91").Value = -3178212.8
$ws.Range("H96").Value = 1769
$ws.Range("I96").Value = 1692
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 5076
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -3703
$ws.Range("N96").Value = -8746
$ws.Range("H98").Value = 3049.75
$ws.Range("I98").Value = 2562.8125
$ws.Range("J98").Value = 4997.5
$ws.Range("K98").Value = 2562.8125
$ws.Range("L98").Value = 4997.5
$ws.Range("M98").Value = -1064.8125
$ws.Range("N98").Value = -7993.5
$ws.Range("H113").Value = 2516.2222
$ws.Range("I113").Value = 2506.6316
$ws.Range("J113").Value = 2539
$ws.Range("K113").Value = 2506.6316
$ws.Range("L113").Value = 2539
$ws.Range("M113").Value = 747.3683999999998
$ws.Range("N113").Value = -9047
$ws.Range("H116").Value = 3467.7
$ws.Range("I116").Value = 3237.8333
$ws.Range("J116").Value = 3812.5
$ws.Range("K116").Value = 3237.8333
$ws.Range("L116").Value = 3812.5
$ws.Range("M116").Value = 204.1667000000002
$ws.Range("N116").Value = -10696.5
$ws.Range("H122").Value = 3049.75
$ws.Range("I122").Value = 2562.8125
$ws.Range("J122").Value = 4997.5
$ws.Range("K122").Value = 7688.4375
$ws.Range("L122").Value = 14992.5
$ws.Range("M122").Value = -5238.4375
$ws.Range("N122").Value = -19892.5
$ws.Range("H138").Value = 1413.66
$ws.Range("I138").Value = 674.30304
$ws.Range("J138").Value = 2848.8823
$ws.Range("K138").Value = 2022.90912
$ws.Range("L138").Value = 8546.6469
$ws.Range("M138").Value = 3117.09088
$ws.Range("N138").Value = -18826.6469
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18167.348
$ws.Range("I32").Value = 20677.484
$ws.Range("K32").Value = 20677.484
$ws.Range("M32").Value = -20390.484
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4401
$ws.Range("I105").Value = 4089.2942
$ws.Range("J105").Value = 4842.5835
$ws.Range("K105").Value = 4089.2942
$ws.Range("L105").Value = 4842.5835
$ws.Range("M105").Value = -2342.2942
$ws.Range("N105").Value = -8336.583500000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2349.3276
$ws.Range("I31").Value = 2468.4443
$ws.Range("J31").Value = 2154.4092
$ws.Range("K31").Value = 2468.4443
$ws.Range("L31").Value = 2154.4092
$ws.Range("M31").Value = -2173.4443
$ws.Range("N31").Value = -2744.4092
$ws.Range("H34").Value = 2349.3276
$ws.Range("I34").Value = 2468.4443
$ws.Range("J34").Value = 2154.4092
$ws.Range("K34").Value = 2468.4443
$ws.Range("L34").Value = 2154.4092
$ws.Range("M34").Value = -2266.4443
$ws.Range("N34").Value = -2558.4092
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 850.6429000000001
$ws.Range("I5").Value = 882.63635
$ws.Range("K5").Value = 2647.90905
$ws.Range("M5").Value = -2535.90905
$ws.Range("H122").Value = 286142
$ws.Range("I122").Value = 188
$ws.Range("K122").Value = 1692
$ws.Range("M122").Value = 758
$ws.Range("H135").Value = 850.6429000000001
$ws.Range("I135").Value = 882.63635
$ws.Range("K135").Value = 7943.72715
$ws.Range("M135").Value = -5408.72715
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13608100
$ws.Range("I11").Value = 17004500
$ws.Range("J11").Value = 22500
$ws.Range("K11").Value = 17004500
$ws.Range("L11").Value = 22500
$ws.Range("M11").Value = -17004361
$ws.Range("N11").Value = -22778
$ws.Range("H70").Value = 4280.769
$ws.Range("I70").Value = 4088.2354
$ws.Range("J70").Value = 4644.4443
$ws.Range("K70").Value = 4088.2354
$ws.Range("L70").Value = 4644.4443
$ws.Range("M70").Value = -3818.2354
$ws.Range("N70").Value = -5184.4443
$ws.Range("H73").Value = 4280.769
$ws.Range("I73").Value = 4088.2354
$ws.Range("J73").Value = 4644.4443
$ws.Range("K73").Value = 4088.2354
$ws.Range("L73").Value = 4644.4443
$ws.Range("M73").Value = -3152.2354
$ws.Range("N73").Value = -6516.4443
$ws.Range("H97").Value = 1207.0834
$ws.Range("I97").Value = 1027
$ws.Range("J97").Value = 1891.4
$ws.Range("K97").Value = 1027
$ws.Range("L97").Value = 1891.4
$ws.Range("M97").Value = -531
$ws.Range("N97").Value = -2883.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3706326.5
$ws.Range("I7").Value = 2758.45
$ws.Range("J7").Value = 11113463
$ws.Range("K7").Value = 2758.45
$ws.Range("L7").Value = 11113463
$ws.Range("M7").Value = -2646.45
$ws.Range("N7").Value = -11113687
$ws.Range("H61").Value = 5256.231
$ws.Range("I61").Value = 5498.4165
$ws.Range("J61").Value = 2350
$ws.Range("K61").Value = 5498.4165
$ws.Range("L61").Value = 2350
$ws.Range("M61").Value = -5296.4165
$ws.Range("N61").Value = -2754
$ws.Range("H113").Value = 5256.231
$ws.Range("I113").Value = 5498.4165
$ws.Range("J113").Value = 2350
$ws.Range("K113").Value = 5498.4165
$ws.Range("L113").Value = 2350
$ws.Range("M113").Value = -3328.4165
$ws.Range("N113").Value = -6690
$ws.Range("H122").Value = 7972.4707
$ws.Range("I122").Value = 11057.454
$ws.Range("J122").Value = 2316.6667
$ws.Range("K122").Value = 33172.362
$ws.Range("L122").Value = 6950.000100000001
$ws.Range("M122").Value = -30722.362
$ws.Range("N122").Value = -11850.0001
$ws.Range("H126").Value = 3706326.5
$ws.Range("I126").Value = 2758.45
$ws.Range("J126").Value = 11113463
$ws.Range("K126").Value = 8275.349999999999
$ws.Range("L126").Value = 33340389
$ws.Range("M126").Value = -5805.349999999999
$ws.Range("N126").Value = -33345329
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 62503.75
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H132").Value = 637.0454999999999
$ws.Range("I132").Value = 527
$ws.Range("J132").Value = 1187.2727
$ws.Range("K132").Value = 1581
$ws.Range("L132").Value = 3561.8181
$ws.Range("M132").Value = 949
$ws.Range("N132").Value = -8621.8181
$ws.Range("H136").Value = 505.41666
$ws.Range("I136").Value = 283.16666
$ws.Range("J136").Value = 1616.6666
$ws.Range("K136").Value = 849.4999799999999
$ws.Range("L136").Value = 4849.9998
$ws.Range("M136").Value = 1700.50002
$ws.Range("N136").Value = -9949.9998
